$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Thermostat) - shorten comment
$ws.Range("E7").Value = "From PC37"

# Row 11 (Join plate) - reword comment
$ws.Range("E11").Value = "Attach the oil sump shell to the engine"

# Row 12 (Anti-planing plate) - reword comment
$ws.Range("E12").Value = "Limit oil displacement in the oil sump"

# Row 15 (Shifter axis) - reword comment
$ws.Range("E15").Value = "Shaft between Shifter and gear motor"

# Row 16 (Shifter gear) - remove comment (now empty) and shrink row height
$ws.Range("E16").Value = ""
$ws.Rows.Item(16).RowHeight = 16.8

# Update view/selection state
$ws.Range("I3:I19").Select()
